$wb = $excel.ActiveWorkbook

# The workbook has duplicated data across the "展览" and "全部类型" sheets.
# Both need the same "想去人数" (F column) updates applied.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 9223
    $ws.Range("F4").Value = 487
}
